# Update "想去人数" (interest count) values in column F on the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# Row -> new value map (column F) shared by both sheets.
$updates = @{
    2  = 8345
    3  = 7789
    4  = 123
    5  = 190
    9  = 118
    10 = 164
    12 = 708
    13 = 128
    14 = 1336
    15 = 62
    19 = 123
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
